$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Prey Class" table (Table2) gains a new row for the lobster prey class.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "lobster"
$ws.Range("C16").Value = "Lobster (southern CA only)"

# Match the author's final selection (cell B16, not C16).
$ws.Range("B16").Select()
